$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Score summary block (rows 10-12) ---
# Copy the "mtitleStyle" formatting (already used by A9) onto the row-label
# cells in column A for rows 10-12.
$ws.Range("A9").Copy()
$ws.Range("A10:A12").PasteSpecial(-4122)  # xlPasteFormats

# Right / Not Attempt / Max counts
$ws.Range("B10").Value = 8
$ws.Range("D10").Value = 20
$ws.Range("E10").Value = 28

# Marking scheme: +4 per right answer, -1 penalty (now numeric, not text)
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# Total score
$ws.Range("B12").Value = 32
$ws.Range("E12").Value = "32/112"

# --- Remove the second and third "Student Ans / Correct Ans" blocks ---
$ws.Range("G15:H40").Clear()
$ws.Range("D19:E40").Clear()

# --- Fill in the student's answers that were correct (column A) ---
# Copy formatting (the green "correctStyle") from B10, which already uses it.
$ws.Range("B10").Copy()

$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("A22").Value = "Option D"

$ws.Range("A24").PasteSpecial(-4122)
$ws.Range("A24").Value = "Option A"

$ws.Range("A27").PasteSpecial(-4122)
$ws.Range("A27").Value = "Option A"

$ws.Range("A28").PasteSpecial(-4122)
$ws.Range("A28").Value = "Option D"

$ws.Range("A29").PasteSpecial(-4122)
$ws.Range("A29").Value = "Option D"

$ws.Range("A34").PasteSpecial(-4122)
$ws.Range("A34").Value = "Option B"

$ws.Range("A35").PasteSpecial(-4122)
$ws.Range("A35").Value = "Option D"

$ws.Range("A37").PasteSpecial(-4122)
$ws.Range("A37").Value = "Option A"
